$d = $word.ActiveDocument

$replacements = @(
    @("93÷2=", "47÷6="),
    @("67÷8=", "59÷4="),
    @("60÷9=", "76÷8="),
    @("75÷7=", "23÷4="),
    @("53÷8=", "95÷6="),
    @("59÷2=", "39÷3="),
    @("39÷4=", "25÷6="),
    @("45÷8=", "56÷7="),
    @("87÷8=", "55÷3="),
    @("71÷9=", "62÷8="),
    @("23÷5=", "53÷4="),
    @("82÷8=", "39÷8="),
    @("24÷4=", "42÷6="),
    @("89÷2=", "93÷5="),
    @("79÷3=", "66÷4="),
    @("97÷3=", "84÷7="),
    @("90÷8=", "71÷8="),
    @("60÷3=", "60÷5="),
    @("91÷3=", "11÷7="),
    @("80÷7=", "69÷5="),
    @("48÷4=", "20÷3="),
    @("62÷4=", "64÷4="),
    @("85÷6=", "10÷2="),
    @("65÷9=", "89÷8="),
    @("21÷6=", "73÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
